$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.652.34"
$ws.Range("E2").Value = "  +1.05%  "
$ws.Range("D3").Value = "1.636.21"
$ws.Range("E3").Value = "  +1.60%  "
$ws.Range("E4").Value = "  +0.00%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "213.12"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.01%  "
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("E7").Value = "  +1.00%  "
$ws.Range("E9").Value = "  +1.02%  "
$ws.Range("E10").Value = "  +3.45%  "
$ws.Range("D12").Value = "1.864.32"
$ws.Range("E12").Value = "  +1.52%  "
$ws.Range("D13").Value = "1.650.43"
$ws.Range("E13").Value = "  +2.52%  "
$ws.Range("E14").Value = "  +0.44%  "
$ws.Range("E15").Value = "  +1.30%  "
$ws.Range("D16").Value = "26.649.73"
$ws.Range("E16").Value = "  +1.05%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "63.00"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +1.44%  "
$ws.Range("E18").Value = "  +1.49%  "
$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +0.02%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "210.07"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +3.31%  "
$ws.Range("E21").Value = "  +0.62%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "9.39"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.56%  "
$ws.Range("E23").Value = "  +1.60%  "
$ws.Range("E24").Value = "  +3.15%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "146.04"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +0.84%  "
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("E27").Value = "  -1.25%  "
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "6.70"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +1.79%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "15.40"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +1.11%  "
$ws.Range("E30").Value = "  +4.77%  "
$ws.Range("E31").Value = "  -0.46%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "3.23"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +0.73%  "
$ws.Range("E33").Value = "  +0.76%  "
$ws.Range("E34").Value = "  +0.70%  "
$ws.Range("E35").Value = "  -0.60%  "
$ws.Range("D36").Value = "1.167.34"
$ws.Range("E36").Value = "  -0.03%  "
$ws.Range("E37").Value = "  +0.16%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.809"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +2.38%  "
$ws.Range("E40").Value = "  -0.26%  "
$ws.Range("E41").Value = "  +0.08%  "
$ws.Range("E42").Value = "  +1.50%  "
$ws.Range("E43").Value = "  +1.77%  "
$ws.Range("D44").Value = "1.773.96"
$ws.Range("E44").Value = "  +1.30%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "92.17"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +0.20%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "1.56"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +1.56%  "
$ws.Range("E47").Value = "  +7.25%  "
$ws.Range("E48").Value = "  +0.52%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.0512"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +0.91%  "
$ws.Range("E50").Value = "  +0.46%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "7.52"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +3.77%  "
